# lines_states.xlsx: add contingency lines line7/line8 and refresh extr1..extr8
# (commit: "contingencies with rene fine")
#
# Final layout (rows 2-7 / line1-line6 stay exactly as-is):
#   row 8  -> line7  (new)
#   row 9  -> line8  (new)
#   row 10 -> extr1  (was row 8)
#   row 11 -> extr2  (was row 9)
#   row 12 -> extr3  (was row 10)
#   row 13 -> extr4  (was row 11)
#   row 14 -> extr5  (was row 12)
#   row 15 -> extr6  (was row 13)
#   row 16 -> extr7  (new row, was row 14)
#   row 17 -> extr8  (new row, was row 15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two brand-new rows (16 and 17) the same formatting as the other
# "id" cells in column A (bold / bordered / centered style used by A2:A15).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

# Shift the extr* rows down two positions, working from the bottom up so we
# never overwrite data we still need to read.

# extr8: old row 15 -> new row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# extr7: old row 14 -> new row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# extr6: old row 13 -> new row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# extr5: old row 12 -> new row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# extr4: old row 11 -> new row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# extr3: old row 10 -> new row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# extr2: old row 9 -> new row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# extr1: old row 8 -> new row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# New row: line7 -> row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row: line8 -> row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true
